$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Updated G column values (Tempo)
$ws.Range("G2").Value = 0.01861023902893066
$ws.Range("G3").Value = 0.01877784729003906
$ws.Range("G4").Value = 0.02029132843017578

# New H column values (Tempo Heuristica)
$ws.Range("H2").Value = 0.00385284423828125
$ws.Range("H3").Value = 0.003971099853515625
$ws.Range("H4").Value = 0.00370478630065918

# New I column values (Tempo Total)
$ws.Range("I2").Value = 0.02246308326721191
$ws.Range("I3").Value = 0.02274894714355469
$ws.Range("I4").Value = 0.02399611473083496
